$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the B1/C1 header labels, dropping the redundant "@ Nirvi" suffix:
#   B1 used to say "... [Nirvi] @ Nirvi", now becomes the plain "[Nirvi]" label
#   C1 used to say "... [учебная] @ Nirvi | ...", now becomes the "учебная" label (no "@ Nirvi")
# (C1 is updated first so the shared-string table ends up ordered the same way
# it is in the target workbook.)
$ws.Range("C1").Value() = "Ижорский (сойкинский) [учебная] | в учебной системе"
$ws.Range("B1").Value() = "Ижорский (сойкинский) [Nirvi]"

# Header row is shorter now that the labels dropped the "@ Nirvi" suffix
$ws.Rows(1).RowHeight() = 30

# Selection moves from D1 to B2
[void]$ws.Range("B2").Select()
